# Helper: write a text value to a cell without Excel auto-typing it as a
# number (this matters for things like "30%" or "1.0" which COM would
# otherwise convert into a percentage/number). We temporarily force a text
# number format, assign the value, then restore the normal cell formatting
# by pasting formats only from a cell that already carries the correct
# style, so the final style index matches the other data cells.
function Set-TextValue($ws, $excel, $targetAddr, $cleanFormatAddr, $text) {
    $ws.Range($targetAddr).NumberFormat = "@"
    $ws.Range($targetAddr).Value = $text
    $ws.Range($cleanFormatAddr).Copy()
    $ws.Range($targetAddr).PasteSpecial(-4122)
    $excel.CutCopyMode = 0
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Insert a new row at position 6 -------------------------------------
# This shifts the old "TOTAL" row (row 6) down to row 7 and Excel
# automatically updates the sheet dimension and the merged-cell range.
$ws.Rows("6:6").Insert()

# Copy the formatting of the data row above (row 5) onto the freshly
# inserted row 6 so its cell style matches the other data rows instead of
# the default style Excel assigns to a blank inserted row.
$ws.Range("A5:G5").Copy()
$ws.Range("A6:G6").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- 2. Row 2 : IMITANCIOMETRIA --------------------------------------------
$ws.Range("C2").Value = "R$ 69,00"
$ws.Range("G2").Value = "R$ 20,70"

# --- 3. Row 3 : now LOGOAUDIOMETRIA (LDV-IRF-LRF) ---------------------------
$ws.Range("A3").Value = "LOGOAUDIOMETRIA (LDV-IRF-LRF)"
$ws.Range("C3").Value = "R$ 78,75"
$ws.Range("G3").Value = "R$ 23,62"

# --- 4. Row 4 : now AUDIOMETRIA EM CAMPO LIVRE ------------------------------
$ws.Range("A4").Value = "AUDIOMETRIA EM CAMPO LIVRE"
$ws.Range("C4").Value = "R$ 60,39"
$ws.Range("G4").Value = "R$ 18,12"

# --- 5. Row 5 : now AUDIOMETRIA TONAL LIMIAR (VIA AEREA / OSSEA) -----------
$ws.Range("A5").Value = "AUDIOMETRIA TONAL LIMIAR (VIA AEREA / OSSEA)"
$ws.Range("C5").Value = "R$ 63,00"
$ws.Range("D5").Value = "BÁSICO"
$ws.Range("E5").Value = "SIM"
Set-TextValue $ws $excel "F5" "F2" "30%"
$ws.Range("G5").Value = "R$ 18,90"

# --- 6. Row 6 : new row - ACOMPANHAMENTO ... --------------------------------
$ws.Range("A6").Value = "ACOMPANHAMENTO DE PACIENTE P/ ADAPTACAO DE APARELHO DE AMPLIFICACAO SONORA INDIVIDUAL (AASI) UNI / BILATERAL"
Set-TextValue $ws $excel "B6" "B2" "1.0"
$ws.Range("C6").Value = "R$ 32,52"
$ws.Range("D6").Value = "CONSULTA"
$ws.Range("E6").Value = "SIM"
Set-TextValue $ws $excel "F6" "F2" "30%"
$ws.Range("G6").Value = "R$ 9,76"

# --- 7. Row 7 : TOTAL (shifted down from old row 6) -------------------------
$ws.Range("G7").Value = "R$ 91,10"
